# The document has two Pearson-logo pictures (in the "default" and
# "first page" footers) and one BTec-logo picture (in the "first page"
# header). Each picture's drawing object carries a name (surfaced on
# both <wp:docPr> and the nested <pic:cNvPr>) that was swapped with its
# sibling image's name in the commit this script reproduces:
#
#   footer (default)    : id=2  image1.png -> image2.png
#   footer (first page)  : id=3  image1.png -> image2.png
#   header (first page)  : id=1  image2.jpg -> image1.jpg
#
# InlineShape objects are anchored inline in the header/footer ranges,
# so we reach them through Sections(1).Headers/Footers rather than
# ActiveDocument.InlineShapes (which only sees the main body story).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer, default (primary) story -> word/footer2.xml, docPr id="2" ---
$footerDefault = $sec.Footers.Item(1)
if ($footerDefault.Exists -and $footerDefault.Range.InlineShapes.Count -ge 1) {
    $pearsonDefault = $footerDefault.Range.InlineShapes.Item(1)
    $pearsonDefault.Name = "image2.png"
}

# --- Footer, first-page story -> word/footer1.xml, docPr id="3" ---
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $pearsonFirst = $footerFirst.Range.InlineShapes.Item(1)
    $pearsonFirst.Name = "image2.png"
}

# --- Header, first-page story -> word/header1.xml, docPr id="1" ---
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $btec = $headerFirst.Range.InlineShapes.Item(1)
    $btec.Name = "image1.jpg"
}
